# KIBON-2653 add other columns again
$wb = $excel.ActiveWorkbook

# --- Gemeinden sheet: insert a new column N (14) ---
$wsGemeinden = $wb.Worksheets.Item("Gemeinden")
$wsGemeinden.Columns.Item(14).Insert()

$wsGemeinden.Cells.Item(8, 14).Value = "{betreuungsstundenFaktor3}"
$wsGemeinden.Cells.Item(7, 14).Value = "1_Betreuungsstunden_Faktor_3"

# --- Tagesschulen sheet: swap the Kinder_Faktor_3 / Kinder_Faktor_15 columns (I and J) back ---
$wsTagesschulen = $wb.Worksheets.Item("Tagesschulen")

$i1 = $wsTagesschulen.Cells.Item(1, 9).Value2
$j1 = $wsTagesschulen.Cells.Item(1, 10).Value2
$wsTagesschulen.Cells.Item(1, 9).Value = $j1
$wsTagesschulen.Cells.Item(1, 10).Value = $i1

$i2 = $wsTagesschulen.Cells.Item(2, 9).Value2
$j2 = $wsTagesschulen.Cells.Item(2, 10).Value2
$wsTagesschulen.Cells.Item(2, 9).Value = $j2
$wsTagesschulen.Cells.Item(2, 10).Value = $i2
